$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise row 7 (2025-07-23 / BEMOL S/A / 387489 / POWER BANK 5000mAh PN-952):
#     estoque_atualizado and desvio_padrao were recomputed ---
$ws.Cells.Item(7, 7).Value = -245
$ws.Cells.Item(7, 9).Value = 0.16

# id_venda (col A) and id_produto (col D) are stored as text in this sheet, so
# force a text number format before typing the new values to keep Excel from
# auto-coercing them into dates / numbers, then drop back to the default
# "Normal" style (matches how the rest of the column is formatted).
$ws.Range("A9:A10").NumberFormat = "@"
$ws.Range("D9:D10").NumberFormat = "@"

# --- New row 9: 2025-07-29 / BEMOL S/A / 389675 / KIT SMARTWATCH INOVA PULSEIRA PRETO LISA (FONE+FONTE+CABO) ---
$ws.Cells.Item(9, 1).Value = "2025-07-29"
$ws.Cells.Item(9, 2).Value = 2
$ws.Cells.Item(9, 3).Value = "BEMOL S/A"
$ws.Cells.Item(9, 4).Value = "389675"
$ws.Cells.Item(9, 5).Value = 49177
$ws.Cells.Item(9, 6).Value = "KIT SMARTWATCH INOVA PULSEIRA PRETO LISA (FONE+FONTE+CABO)"
$ws.Cells.Item(9, 7).Value = -65
$ws.Cells.Item(9, 8).Value = 1.03
$ws.Cells.Item(9, 9).Value = 0.18

# --- New row 10: 2025-07-29 / BEMOL S/A / 389699 / FONE DE OUVIDO TIPO-C EJ-105 ---
$ws.Cells.Item(10, 1).Value = "2025-07-29"
$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(10, 3).Value = "BEMOL S/A"
$ws.Cells.Item(10, 4).Value = "389699"
$ws.Cells.Item(10, 5).Value = 48696
$ws.Cells.Item(10, 6).Value = "FONE DE OUVIDO TIPO-C EJ-105"
$ws.Cells.Item(10, 7).Value = -73
$ws.Cells.Item(10, 8).Value = 1.06
$ws.Cells.Item(10, 9).Value = 0.23

$ws.Range("A9:A10").Style = "Normal"
$ws.Range("D9:D10").Style = "Normal"
